$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Periodo Mora" (column E) history for trabajador YERALDIN PAOLA MERCADO OLIVERO
# (rows 16-28) is reversed into descending order, and the trabajador
# MARIELE DEL MAR PIÑA PAJARO record (previously row 27) moves down to row 29.

$ws.Range("E16").Value = "2104"
$ws.Range("F16").Value = 30430

$ws.Range("E17").Value = "2103"

$ws.Range("E18").Value = "2102"

$ws.Range("E19").Value = "2101"

$ws.Range("E20").Value = "2012"

$ws.Range("E21").Value = "2011"

# Row 22 (2010) is unchanged.

$ws.Range("E23").Value = "2009"

$ws.Range("E24").Value = "2008"

$ws.Range("E25").Value = "2007"

$ws.Range("E26").Value = "2006"

# Row 27 previously held MARIELE DEL MAR PIÑA PAJARO's record; it now holds
# another YERALDIN period (2005).
$ws.Range("C27").Value = "1047438450"
$ws.Range("D27").Value = "YERALDIN PAOLA MERCADO OLIVERO"
$ws.Range("E27").Value = "2005"
$ws.Range("F27").Value = 35112
$ws.Range("G27").Value = 1200000

$ws.Range("E28").Value = "2004"

# Row 29 now holds MARIELE DEL MAR PIÑA PAJARO's record (moved from row 27).
$ws.Range("C29").Value = "1143401993"
$ws.Range("D29").Value = "MARIELE DEL MAR PIÑA PAJARO"
$ws.Range("E29").Value = "2103"
$ws.Range("F29").Value = 9691
$ws.Range("G29").Value = 908526
